# Lecture partielle de l'EDT M1 MIAGE.
# Shift the schedule dates forward by 3 years (1096 days) and update the
# corresponding French weekday labels accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 : vendredi 13/01/2023 -> mardi 13/01/2026
$ws.Range("A2").Value = 46035.0
$ws.Range("B2").Value = "mardi"

# Row 4 : dimanche 15/01/2023 -> jeudi 15/01/2026
$ws.Range("A4").Value = 46037.0
$ws.Range("B4").Value = "jeudi"

# Row 7 : vendredi 20/01/2023 -> mardi 20/01/2026
$ws.Range("A7").Value = 46042.0
$ws.Range("B7").Value = "mardi"

# Row 9 : dimanche 22/01/2023 -> jeudi 22/01/2026
$ws.Range("A9").Value = 46044.0
$ws.Range("B9").Value = "jeudi"

# Row 12 : vendredi 12/05/2023 -> mardi 12/05/2026
$ws.Range("A12").Value = 46154.0
$ws.Range("B12").Value = "mardi"

# Row 15 : vendredi 19/05/2023 -> mardi 19/05/2026
$ws.Range("A15").Value = 46161.0
$ws.Range("B15").Value = "mardi"
